$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.231.78"
$ws.Range("E2").Value = "  +2.16%  "
$ws.Range("D3").Value = "3.365.95"
$ws.Range("E3").Value = "  +1.72%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Value = "'588.81"
$ws.Range("E5").Value = "  +6.36%  "
$ws.Range("D6").Value = "'188.34"
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").Value = "'0.600"
$ws.Range("E8").Value = "  +3.57%  "
$ws.Range("D9").Value = "'0.184"
$ws.Range("E9").Value = "  +2.92%  "
$ws.Range("E10").Value = "  +1.19%  "
$ws.Range("D11").Value = "'47.69"
$ws.Range("E11").Value = "  +2.17%  "
$ws.Range("D12").Value = "'0.0000275"
$ws.Range("E12").Value = "  +3.19%  "
$ws.Range("D13").Value = "'662.33"
$ws.Range("E13").Value = "  +11.00%  "
$ws.Range("D14").Value = "3.901.87"
$ws.Range("E14").Value = "  +1.71%  "
$ws.Range("D15").Value = "'8.64"
$ws.Range("E15").Value = "  +0.46%  "
$ws.Range("D16").Value = "67.247.88"
$ws.Range("E16").Value = "  +2.15%  "
$ws.Range("E17").Value = "  +1.06%  "
$ws.Range("D18").Value = "'18.07"
$ws.Range("E18").Value = "  +1.12%  "
$ws.Range("D19").Value = "3.360.55"
$ws.Range("E19").Value = "  +1.75%  "
$ws.Range("E20").Value = "  +2.07%  "
$ws.Range("D21").Value = "'0.911"
$ws.Range("E21").Value = "  +1.60%  "
$ws.Range("D22").Value = "'18.15"
$ws.Range("E22").Value = "  -1.52%  "
$ws.Range("E23").Value = "  +0.80%  "
$ws.Range("D24").Value = "'101.26"
$ws.Range("E24").Value = "  +0.84%  "
$ws.Range("E25").Value = "  +2.39%  "
$ws.Range("D26").Value = "'2.85"
$ws.Range("E26").Value = "  +4.64%  "
$ws.Range("D27").Value = "'9.79"
$ws.Range("E27").Value = "  +3.54%  "
$ws.Range("D28").Value = "'32.27"
$ws.Range("D29").Value = "'8.74"
$ws.Range("E29").Value = "  +0.86%  "
$ws.Range("D30").Value = "'6.88"
$ws.Range("E30").Value = "  +2.69%  "
$ws.Range("D31").Value = "'614.29"
$ws.Range("E31").Value = "  +7.71%  "
$ws.Range("D32").Value = "'3.92"
$ws.Range("E32").Value = "  +2.73%  "
$ws.Range("E34").Value = "  +3.03%  "
$ws.Range("D35").Value = "3.883.42"
$ws.Range("E35").Value = "  +4.81%  "
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("D37").Value = "'55.43"
$ws.Range("E37").Value = "  -2.30%  "
$ws.Range("D38").Value = "'2.84"
$ws.Range("E38").Value = "  +7.49%  "
$ws.Range("E39").Value = "  +2.78%  "
$ws.Range("D40").Value = "'33.87"
$ws.Range("E40").Value = "  +0.82%  "
$ws.Range("D41").Value = "'3.27"
$ws.Range("E41").Value = "  +1.48%  "
$ws.Range("D42").Value = "0.0₃0707"
$ws.Range("E42").Value = "  +1.73%  "
$ws.Range("E43").Value = "  +2.96%  "
$ws.Range("E44").Value = "  -0.96%  "
$ws.Range("D45").Value = "'0.0424"
$ws.Range("E45").Value = "  +2.01%  "
$ws.Range("E46").Value = "  +1.35%  "
$ws.Range("D47").Value = "'2.60"
$ws.Range("E47").Value = "  +2.14%  "
$ws.Range("D48").Value = "'1.00"
$ws.Range("E48").Value = "  +0.32%  "
$ws.Range("D49").Value = "'2.89"
$ws.Range("E49").Value = "  -17.39%  "
$ws.Range("D50").Value = "'1.35"
$ws.Range("E50").Value = "  +9.25%  "
$ws.Range("D51").Value = "'129.69"
$ws.Range("E51").Value = "  +4.97%  "
